# Update "想去人数" (number of people interested) counts in the
# "展览" and "全部类型" sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1168
    3  = 592
    5  = 34
    6  = 165
    10 = 5406
    11 = 4857
    13 = 40
    15 = 51
    16 = 194
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
